{"js": "// Ordered list of [oldText, newText] pairs, one per non-empty paragraph,\n// in document order (the title line, then each filled table cell).\nconst replacements = [\n  [\"2024-09-05 Thursday\", \"2024-09-06 Friday\"],\n  [\"70\u00d755=\", \"29\u00d761=\"],\n  [\"87\u00d740=\", \"51\u00d756=\"],\n  [\"16\u00d714=\", \"92\u00d732=\"],\n  [\"96\u00d722=\", \"87\u00d773=\"],\n  [\"64\u00d759=\", \"54\u00d731=\"],\n  [\"21\u00d740=\", \"51\u00d772=\"],\n  [\"74\u00d752=\", \"80\u00d720=\"],\n  [\"41\u00d787=\", \"40\u00d768=\"],\n  [\"15\u00d769=\", \"46\u00d743=\"],\n  [\"80\u00d741=\", \"87\u00d788=\"],\n  [\"19\u00d791=\", \"95\u00d711=\"],\n  [\"53\u00d712=\", \"66\u00d743=\"],\n  [\"53\u00d766=\", \"31\u00d717=\"],\n  [\"70\u00d755=\", \"80\u00d781=\"],\n  [\"90\u00d752=\", \"90\u00d713=\"],\n  [\"96\u00d784=\", \"60\u00d733=\"],\n  [\"36\u00d773=\", \"14\u00d797=\"],\n  [\"40\u00d787=\", \"98\u00d796=\"],\n  [\"18\u00d733=\", \"53\u00d767=\"],\n  [\"13\u00d775=\", \"64\u00d747=\"],\n  [\"93\u00d759=\", \"31\u00d714=\"],\n  [\"90\u00d759=\", \"29\u00d789=\"],\n  [\"31\u00d787=\", \"29\u00d788=\"],\n  [\"81\u00d713=\", \"30\u00d788=\"],\n  [\"61\u00d776=\", \"58\u00d795=\"],\n];\n\n// `body.paragraphs` walks the whole document in order, including the\n// paragraphs that live inside table cells, so the non-empty paragraphs\n// line up 1:1 (in order) with the `replacements` list above.\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nlet repIndex = 0;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (repIndex >= replacements.length) break;\n  const para = paragraphs.items[i];\n  const text = para.text;\n  if (text === \"\") continue; // blank answer cells are left untouched\n\n  const [oldText, newText] = replacements[repIndex];\n  if (text !== oldText) {\n    throw new Error(\n      `Unexpected paragraph text at position ${repIndex}: expected \"${oldText}\" but found \"${text}\"`\n    );\n  }\n  para.insertText(newText, Word.InsertLocation.replace);\n  repIndex++;\n}\n\nawait context.sync();\n\nif (repIndex !== replacements.length) {\n  throw new Error(`Only applied ${repIndex} of ${replacements.length} replacements`);\n}\n", "ps1": "# Ordered list of old/new text pairs. The first entry is the title line\n# (a plain paragraph before the table); the rest are the multiplication\n# expressions, in document order (row by row, left to right). Blank\n# answer cells are skipped, so the filled cells line up 1:1 (in order)\n# with this list.\n$replacements = @(\n    @(\"2024-09-05 Thursday\", \"2024-09-06 Friday\"),\n    @(\"70\u00d755=\", \"29\u00d761=\"),\n    @(\"87\u00d740=\", \"51\u00d756=\"),\n    @(\"16\u00d714=\", \"92\u00d732=\"),\n    @(\"96\u00d722=\", \"87\u00d773=\"),\n    @(\"64\u00d759=\", \"54\u00d731=\"),\n    @(\"21\u00d740=\", \"51\u00d772=\"),\n    @(\"74\u00d752=\", \"80\u00d720=\"),\n    @(\"41\u00d787=\", \"40\u00d768=\"),\n    @(\"15\u00d769=\", \"46\u00d743=\"),\n    @(\"80\u00d741=\", \"87\u00d788=\"),\n    @(\"19\u00d791=\", \"95\u00d711=\"),\n    @(\"53\u00d712=\", \"66\u00d743=\"),\n    @(\"53\u00d766=\", \"31\u00d717=\"),\n    @(\"70\u00d755=\", \"80\u00d781=\"),\n    @(\"90\u00d752=\", \"90\u00d713=\"),\n    @(\"96\u00d784=\", \"60\u00d733=\"),\n    @(\"36\u00d773=\", \"14\u00d797=\"),\n    @(\"40\u00d787=\", \"98\u00d796=\"),\n    @(\"18\u00d733=\", \"53\u00d767=\"),\n    @(\"13\u00d775=\", \"64\u00d747=\"),\n    @(\"93\u00d759=\", \"31\u00d714=\"),\n    @(\"90\u00d759=\", \"29\u00d789=\"),\n    @(\"31\u00d787=\", \"29\u00d788=\"),\n    @(\"81\u00d713=\", \"30\u00d788=\"),\n    @(\"61\u00d776=\", \"58\u00d795=\")\n)\n\n$d = $word.ActiveDocument\n$repIndex = 0\n\n# --- Title paragraph (the single paragraph that precedes the table) ---\n$titlePara = $d.Paragraphs.Item(1)\n$titleRange = $titlePara.Range\n$titleText = $titleRange.Text.TrimEnd([char]13, [char]7)\n$pair = $replacements[$repIndex]\nif ($titleText -ne $pair[0]) {\n    throw \"Unexpected title text: expected [$($pair[0])] but found [$titleText]\"\n}\n$targetRange = $d.Range($titleRange.Start, $titleRange.Start + $titleText.Length)\n$targetRange.Text = $pair[1]\n$repIndex = $repIndex + 1\n\n# --- Table cells, in row-major / document order ---\n$table = $d.Tables.Item(1)\n$done = $false\nforeach ($row in $table.Rows) {\n    if ($done) { break }\n    foreach ($cell in $row.Cells) {\n        if ($repIndex -ge $replacements.Count) { $done = $true; break }\n        $cellRange = $cell.Range\n        $cellText = $cellRange.Text.TrimEnd([char]13, [char]7)\n        if ($cellText -eq \"\") { continue }\n\n        $pair = $replacements[$repIndex]\n        if ($cellText -ne $pair[0]) {\n            throw \"Unexpected cell text at position $repIndex : expected [$($pair[0])] but found [$cellText]\"\n        }\n        $cellRange.Text = $pair[1]\n        $repIndex = $repIndex + 1\n    }\n}\n\nif ($repIndex -ne $replacements.Count) {\n    throw \"Only applied $repIndex of $($replacements.Count) replacements\"\n}\n"}
